$d = $word.ActiveDocument

# 1. Rewrite the "My C++ skills..." paragraph: merge all of its runs into a
#    single run and append the new "iterating on gameplay systems" sentence.
$p1 = $d.Paragraphs.Item(10)
$r1 = $p1.Range
$r1body = $d.Range($r1.Start, $r1.End - 1)
$r1body.Text = ""
$r1new = $d.Range($r1.Start, $r1.Start)
$r1new.Text = "My C++ skills are at their peak through rigorous practice with the use of pointers and a better understanding of Data Structures from the Collision System and Memory Manager that I created. I love delving into 3D Math and am relearning it in a better way, with a heavy focus on understanding it through geometry and visualizing it, for use specifically in games. The Action games that I have worked on and am currently working on have given me experience in bringing the design, engineering, art and animation in them together, and in collaborating with and learning from the people involved in them, as well as iterating on gameplay systems to get them to their best possible form for the game."

# 2. Move the hidden "_GoBack" bookmark from the "Talk about..." paragraph
#    to the very start of the "Combat, weapons..." paragraph.
$bk = $d.Bookmarks.Item("_GoBack")
$bk.Delete()

$p11 = $d.Paragraphs.Item(11)
$bkRange = $d.Range($p11.Range.Start, $p11.Range.Start)
$d.Bookmarks.Add("_GoBack", $bkRange)

# 3. Merge the two runs of the "Talk about..." paragraph (now that the
#    bookmark no longer separates them) into a single run.
$p14 = $d.Paragraphs.Item(14)
$r14 = $p14.Range
$r14body = $d.Range($r14.Start, $r14.End - 1)
$r14body.Text = ""
$r14new = $d.Range($r14.Start, $r14.Start)
$r14new.Text = "Talk about wanting to create power fantasy that they go for here? -"
